# Row 38 (Idaho) failed on this run: the scraper call for this location
# timed out, so the previously-populated numeric columns (B:H) are wiped
# back to empty cells, the "includes Hispanic Black" flag (J) flips to
# FALSE, and the status message records the timeout instead of "Success!".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear B38:H38 back to blank text cells (no value, no special number
# format/style) instead of leaving Excel's default "blank number" cell.
$cols = @("B", "C", "D", "E", "F", "G", "H")
foreach ($col in $cols) {
    $cell = $ws.Range($col + "38")
    $cell.Value = "'"
    $cell.Style = "Normal"
}

# "Pct Includes Hispanic Black" is no longer TRUE for this failed run.
$ws.Range("J38").Value = $false

# Status column records the error that occurred instead of "Success!".
$ws.Range("O38").Value = "An error occurred. ... TimeoutException('', None, None)"
